# Ultimate fix on backend
# Swap the "tp" (C column) and "fp" (D column) values for each data row,
# then recompute precision (G) and fscore (I) from the corrected tp/fp values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 16 }

for ($r = 2; $r -le $lastRow; $r++) {
    $tpOld = $ws.Cells.Item($r, 3).Value2  # column C
    $fpOld = $ws.Cells.Item($r, 4).Value2  # column D

    $tpNew = $fpOld
    $fpNew = $tpOld

    $ws.Cells.Item($r, 3).Value = $tpNew
    $ws.Cells.Item($r, 4).Value = $fpNew

    $recall = $ws.Cells.Item($r, 8).Value2  # column H

    $precision = $tpNew / ($tpNew + $fpNew)
    $ws.Cells.Item($r, 7).Value = $precision  # column G

    $fscore = (2 * $precision * $recall) / ($precision + $recall)
    $ws.Cells.Item($r, 9).Value = $fscore  # column I
}
